# dist_users_ug16.xlsx — "maj coût et users"
# - deactivate (actif=0) the first 4 users (rows 2-5)
# - clear the "lien_maj" (update link) column for every data row
# - bump the "Version" value on the last row from 4.1 to 4.0
# - leave the selection on C8, matching the author's last click

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# actif column (C): turn off for rows 2-5
$ws.Range("C2").Value = 0
$ws.Range("C3").Value = 0
$ws.Range("C4").Value = 0
$ws.Range("C5").Value = 0

# lien_maj column (V): wipe the stored OneDrive link, keep the hyperlink-style formatting
$ws.Range("V2:V9").ClearContents()

# Version column (U) on row 9: was "4.1", should become "4.0" like every other row.
# Assign via a copy from a sibling cell that already holds the literal text "4.0"
# so Excel doesn't re-interpret it as a number (and no new number format/style
# gets minted in the process).
$ws.Range("U2").Copy($ws.Range("U9"))

# Restore the simple single-cell selection the author ended on
$ws.Range("C8").Select()
